# Apply corrected Diebold-Mariano values to P_valores and Estadisticos_DM sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.155735342972489
$wsP.Range("D2").Value = 0.06703962894406534
$wsP.Range("E2").Value = 0.1218441147767853
$wsP.Range("F2").Value = 0.1206751132599995

$wsP.Range("B3").Value = 0.155735342972489
$wsP.Range("D3").Value = 0.6717334214439286
$wsP.Range("E3").Value = 0.6232646629428999
$wsP.Range("F3").Value = 0.8668519093654812

$wsP.Range("B4").Value = 0.06703962894406534
$wsP.Range("C4").Value = 0.6717334214439286
$wsP.Range("E4").Value = 0.3552301627312495
$wsP.Range("F4").Value = 0.5037495377485375

$wsP.Range("B5").Value = 0.1218441147767853
$wsP.Range("C5").Value = 0.6232646629428999
$wsP.Range("D5").Value = 0.3552301627312495
$wsP.Range("F5").Value = 0.6847938334849961

$wsP.Range("B6").Value = 0.1206751132599995
$wsP.Range("C6").Value = 0.8668519093654812
$wsP.Range("D6").Value = 0.5037495377485375
$wsP.Range("E6").Value = 0.6847938334849961

# --- Sheet: Estadisticos_DM ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -1.500351137043517
$wsE.Range("D2").Value = -1.985459475008506
$wsE.Range("E2").Value = -1.646844032069017
$wsE.Range("F2").Value = -1.652496127919068

$wsE.Range("B3").Value = 1.500351137043517
$wsE.Range("D3").Value = -0.4328277101067428
$wsE.Range("E3").Value = 0.5022994526076398
$wsE.Range("F3").Value = 0.1707650398062218

$wsE.Range("B4").Value = 1.985459475008506
$wsE.Range("C4").Value = 0.4328277101067428
$wsE.Range("E4").Value = 0.9561226991034093
$wsE.Range("F4").Value = 0.6862722402122541

$wsE.Range("B5").Value = 1.646844032069017
$wsE.Range("C5").Value = -0.5022994526076398
$wsE.Range("D5").Value = -0.9561226991034093
$wsE.Range("F5").Value = -0.4144934036874958

$wsE.Range("B6").Value = 1.652496127919068
$wsE.Range("C6").Value = -0.1707650398062218
$wsE.Range("D6").Value = -0.6862722402122541
$wsE.Range("E6").Value = 0.4144934036874958
